$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking.com price values (column D) are plain text in the sheet (e.g. "592.42",
# "1.00", "67.451.44"). Writing such strings straight into a cell would make Excel
# auto-convert the number-looking ones into real numbers (losing formatting/precision),
# so each Price cell is temporarily switched to Text format, given its new text value,
# and then restored to the default (Normal) style so no visible/style side effects remain.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.451.44"
$ws.Range("E2").Value = "  -0.46%  "
Set-TextValue $ws.Range("D3") "2.613.83"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws.Range("D5") "592.42"
$ws.Range("E5").Value = "  -1.83%  "
Set-TextValue $ws.Range("D6") "150.83"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.59%  "
Set-TextValue $ws.Range("D9") "2.608.81"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("E11").Value = "  -0.07%  "
Set-TextValue $ws.Range("D12") "5.18"
$ws.Range("E12").Value = "  -1.43%  "
Set-TextValue $ws.Range("D13") "0.344"
$ws.Range("E13").Value = "  -3.31%  "
Set-TextValue $ws.Range("D14") "27.34"
$ws.Range("E14").Value = "  -2.83%  "
Set-TextValue $ws.Range("D15") "3.081.63"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("E16").Value = "  -2.55%  "
Set-TextValue $ws.Range("D17") "67.330.21"
Set-TextValue $ws.Range("D18") "2.612.72"
$ws.Range("E18").Value = "  -0.31%  "
Set-TextValue $ws.Range("D19") "371.51"
$ws.Range("E19").Value = "  +1.64%  "
Set-TextValue $ws.Range("D20") "11.06"
$ws.Range("E20").Value = "  -2.45%  "
Set-TextValue $ws.Range("D21") "7.38"
$ws.Range("E21").Value = "  -3.35%  "
Set-TextValue $ws.Range("D22") "4.29"
$ws.Range("E22").Value = "  -0.60%  "
Set-TextValue $ws.Range("D23") "4.76"
$ws.Range("E23").Value = "  -4.58%  "
$ws.Range("E24").Value = "  -3.70%  "
Set-TextValue $ws.Range("D25") "73.46"
$ws.Range("E25").Value = "  +4.70%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("E28").Value = "  -0.07%  "
Set-TextValue $ws.Range("D29") "581.06"
$ws.Range("E29").Value = "  -0.84%  "
Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  -0.21%  "
Set-TextValue $ws.Range("D31") "0.0₃0988"
$ws.Range("E31").Value = "  -6.05%  "
$ws.Range("E32").Value = "  -5.56%  "
Set-TextValue $ws.Range("D33") "7.68"
$ws.Range("E33").Value = "  -3.54%  "
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  -3.96%  "
$ws.Range("E37").Value = "  -2.82%  "
Set-TextValue $ws.Range("D38") "157.76"
$ws.Range("E38").Value = "  +1.40%  "
Set-TextValue $ws.Range("D39") "19.08"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E41").Value = "  -1.80%  "
Set-TextValue $ws.Range("D42") "5.24"
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("E44").Value = "  +4.27%  "
$ws.Range("E45").Value = "  +0.03%  "
Set-TextValue $ws.Range("D46") "153.52"
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("E50").Value = "  -5.70%  "
Set-TextValue $ws.Range("D51") "21.40"
$ws.Range("E51").Value = "  +1.38%  "
